$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Forces a numeric-looking string to be stored as literal text (matching the
# source data export, which writes everything as inlineStr) and then resets
# the cell style back to Normal so no stray number-format/quote-prefix style
# is left behind on the cell.
function Set-TextValue($Cell, $Text) {
    $Cell.Value = "'" + $Text
    $Cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") '277.28'
Set-TextValue $ws.Range("E2") '0.95%'
Set-TextValue $ws.Range("G2") '1'

# Row 3
Set-TextValue $ws.Range("D3") '27.34'
Set-TextValue $ws.Range("E3") '2.46%'
Set-TextValue $ws.Range("G3") '1'

# Row 4
Set-TextValue $ws.Range("D4") '4.872'
Set-TextValue $ws.Range("E4") '-0.52%'
Set-TextValue $ws.Range("G4") '1'

# Row 5
Set-TextValue $ws.Range("D5") '0.06416'
Set-TextValue $ws.Range("E5") '1.23%'
Set-TextValue $ws.Range("G5") '1'

# Row 6
Set-TextValue $ws.Range("D6") '6.954'
Set-TextValue $ws.Range("E6") '0.43%'
Set-TextValue $ws.Range("G6") '1'

# Row 7
Set-TextValue $ws.Range("D7") '1.194'
Set-TextValue $ws.Range("E7") '-6.92%'
Set-TextValue $ws.Range("G7") '1'

# Row 8
Set-TextValue $ws.Range("D8") '0.8869'
Set-TextValue $ws.Range("E8") '0.52%'
Set-TextValue $ws.Range("G8") '1'

# Row 9
Set-TextValue $ws.Range("D9") '0.1519'
Set-TextValue $ws.Range("E9") '3.18%'
Set-TextValue $ws.Range("G9") '1'

# Row 10
Set-TextValue $ws.Range("D10") '0.05069'
Set-TextValue $ws.Range("E10") '0.97%'
Set-TextValue $ws.Range("G10") '1'

# Row 11
Set-TextValue $ws.Range("D11") '0.07539'
Set-TextValue $ws.Range("E11") '2.33%'
Set-TextValue $ws.Range("G11") '1'

# Row 12
Set-TextValue $ws.Range("D12") '0.02958'
Set-TextValue $ws.Range("E12") '-4.34%'
Set-TextValue $ws.Range("G12") '1'

# Row 13
Set-TextValue $ws.Range("D13") '0.09016'
Set-TextValue $ws.Range("E13") '-0.26%'
Set-TextValue $ws.Range("G13") '1'

# Row 14
Set-TextValue $ws.Range("D14") '0.001578'
Set-TextValue $ws.Range("E14") '0.35%'
Set-TextValue $ws.Range("G14") '1'

# Row 15
Set-TextValue $ws.Range("D15") '0.0006430'
Set-TextValue $ws.Range("E15") '1.91%'
Set-TextValue $ws.Range("G15") '1'

# Row 16
Set-TextValue $ws.Range("D16") '0.006185'
Set-TextValue $ws.Range("E16") '2.73%'
Set-TextValue $ws.Range("G16") '1'

# Row 17
Set-TextValue $ws.Range("D17") '3.472'
Set-TextValue $ws.Range("E17") '-0.12%'
Set-TextValue $ws.Range("G17") '1'

# Row 18
Set-TextValue $ws.Range("D18") '3.307'
Set-TextValue $ws.Range("E18") '-1.24%'
Set-TextValue $ws.Range("G18") '1'

# Row 19
Set-TextValue $ws.Range("E19") '0.00%'
Set-TextValue $ws.Range("G19") '1'

# Row 20
Set-TextValue $ws.Range("G20") '1'

# Row 21
Set-TextValue $ws.Range("D21") '0.1346'
Set-TextValue $ws.Range("E21") '2.39%'
Set-TextValue $ws.Range("G21") '1'

# Row 22
Set-TextValue $ws.Range("D22") '3.909'
Set-TextValue $ws.Range("E22") '-0.06%'
Set-TextValue $ws.Range("G22") '1'

# Row 23
Set-TextValue $ws.Range("D23") '0.04416'
Set-TextValue $ws.Range("E23") '1.01%'
Set-TextValue $ws.Range("G23") '1'

# Row 24
Set-TextValue $ws.Range("D24") '0.001174'
Set-TextValue $ws.Range("E24") '-0.37%'
Set-TextValue $ws.Range("G24") '1'

# Row 25
Set-TextValue $ws.Range("D25") '0.004282'
Set-TextValue $ws.Range("E25") '16.20%'
Set-TextValue $ws.Range("G25") '1'

# Row 26
Set-TextValue $ws.Range("D26") '0.0001200'
Set-TextValue $ws.Range("E26") '0.00%'
Set-TextValue $ws.Range("G26") '1'

# Row 27
Set-TextValue $ws.Range("E27") '-0.01%'
Set-TextValue $ws.Range("G27") '1'

# Row 28
Set-TextValue $ws.Range("G28") '1'

# Row 29
Set-TextValue $ws.Range("G29") '1'

# Row 30
Set-TextValue $ws.Range("G30") '1'

# Row 31
Set-TextValue $ws.Range("G31") '1'

# Row 32
Set-TextValue $ws.Range("G32") '1'

# Row 33
Set-TextValue $ws.Range("G33") '1'

# Row 34
Set-TextValue $ws.Range("G34") '1'

# Row 35
Set-TextValue $ws.Range("G35") '1'

# Row 36
Set-TextValue $ws.Range("G36") '1'

# Row 37
Set-TextValue $ws.Range("G37") '1'

# Row 38
Set-TextValue $ws.Range("G38") '1'

# Row 39
Set-TextValue $ws.Range("G39") '1'

# Row 40
Set-TextValue $ws.Range("D40") '0.04166'
Set-TextValue $ws.Range("E40") '2.02%'
Set-TextValue $ws.Range("G40") '1'

# Row 41
Set-TextValue $ws.Range("D41") '0.006826'
Set-TextValue $ws.Range("E41") '3.53%'
Set-TextValue $ws.Range("G41") '1'

# Row 42
Set-TextValue $ws.Range("D42") '0.1178'
Set-TextValue $ws.Range("G42") '1'

# Row 43
Set-TextValue $ws.Range("D43") '0.002090'
Set-TextValue $ws.Range("E43") '-0.95%'
Set-TextValue $ws.Range("G43") '1'

# Row 44
Set-TextValue $ws.Range("D44") '0.01174'
Set-TextValue $ws.Range("E44") '-3.19%'
Set-TextValue $ws.Range("G44") '1'

# Row 45
Set-TextValue $ws.Range("D45") '0.00005177'
Set-TextValue $ws.Range("E45") '-2.66%'
Set-TextValue $ws.Range("G45") '1'

# Row 46
$ws.Range("B46").Value = 'CoinbaseStockToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
Set-TextValue $ws.Range("D46") '0.02000'
Set-TextValue $ws.Range("E46") '-22.92%'
Set-TextValue $ws.Range("G46") '1'

# Row 47
$ws.Range("B47").Value = 'BOLO'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
Set-TextValue $ws.Range("D47") '1.487'
Set-TextValue $ws.Range("E47") '-36.90%'
Set-TextValue $ws.Range("G47") '1'

# Row 48
Set-TextValue $ws.Range("G48") '1'

# Row 49
Set-TextValue $ws.Range("G49") '1'

# Row 50
Set-TextValue $ws.Range("G50") '1'

# Row 51
Set-TextValue $ws.Range("G51") '1'
